$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17, shifting existing rows 17:111 down to 18:112
$ws.Rows.Item(17).Insert()

# Populate the new row 17 with the latest weekly price observation
$ws.Range("A17").Value = 5
$ws.Range("B17").Value = "Macroferia Regional de Talca"
$ws.Range("C17").Value = "Maule"
$ws.Range("D17").Value = 44462
$ws.Range("E17").Value = 7
$ws.Range("F17").Value = 100112017
$ws.Range("G17").Value = "Apio"
$ws.Range("H17").Value = "Americana (o)"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 500
$ws.Range("K17").Value = 8500
$ws.Range("L17").Value = 8500
$ws.Range("M17").Value = 8500
$ws.Range("N17").Value = "`$/docena de matas"
$ws.Range("O17").Value = "Provincia del Elquí"
$ws.Range("P17").Value = 1417
$ws.Range("Q17").Value = 6
$ws.Range("R17").Value = "Hortaliza"
